# The "Recorded By" column (column G) on the "Session Analysis Results"
# sheet lists the users/systems that recorded each attendance session.
# Previously some cells read "System, dnasr281@gmail.com"; the order of
# the two names in that list should be swapped to
# "dnasr281@gmail.com, System" everywhere it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
